$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 712, pushing existing rows 712:753 down to 713:754
$ws.Rows.Item(712).Insert()

# Column A holds a date-looking value but is stored as plain text in this
# sheet (inlineStr), so force text formatting before assigning it to avoid
# Excel auto-converting it into a real date serial number, then drop the
# number-format override so the cell ends up unstyled like its neighbours.
$ws.Cells.Item(712, 1).NumberFormat = "@"
$ws.Cells.Item(712, 1).Value = "2026/01/28"
$ws.Cells.Item(712, 1).ClearFormats()

$ws.Cells.Item(712, 2).Value = "水"
$ws.Cells.Item(712, 3).Value = 19
$ws.Cells.Item(712, 4).Value = 53
